# Edit script: transform Review_177 (TeacherLM paper review) into
# Review_176 (Large Language Models as Generalizable Policies for
# Embodied Tasks paper review).
$d = $word.ActiveDocument

# 1) Title heading (paragraph 1, style Heading1)
$d.Paragraphs.Item(1).Range.Text = "Review 176: Large Language Models as Generalizable Policies for Embodied Tasks"

# 2) Bold "Paper: <arxiv link>" line (paragraph 2)
$d.Paragraphs.Item(2).Range.Text = "Paper: https://arxiv.org/abs/2310.17722v2"

# 3) Huggingface papers link (paragraph 4) - new text has a leading
#    space before "https".
$d.Paragraphs.Item(4).Range.Text = " https://huggingface.co/papers/2310.17722"

# 4) Delete the extra blank "Normal" paragraph that sat between the
#    huggingface-link paragraph and the first long Hebrew body
#    paragraph (this was paragraph 6 out of the original 14
#    paragraphs). After this delete, every later paragraph's index
#    shifts down by one, leaving 13 paragraphs total.
$d.Paragraphs.Item(6).Range.Delete()

# 5) Replace the first long Hebrew body paragraph (now paragraph 6,
#    formerly paragraph 7). Its run carries a leftover
#    xml:space="preserve" from the old text (which ended in
#    trailing spaces); the new text has no leading/trailing
#    whitespace, so rather than editing the run in place (which
#    would keep the stale attribute) we insert a brand-new paragraph
#    after paragraph 5 and delete the old one.
$d.Paragraphs.Item(5).Range.InsertParagraphAfter()
$d.Paragraphs.Item(6).Range.Text = "על למידה עם חיזוקים (reinforcement learning) שמעתם כבר? על מודלי שפה בטח שמעתם, נכון? אז היום אנחנו נדבר על השידוך ביניהם.  אזכיר ש-RL היא למעשה משפחת שיטות המאפשרות לאמן מודל on-the-fly. כלומר תוך כדי אימון המודל ניתן ליצור דאטה כל פעם שהמודל מתאמן ולהמשיך לאמן עליו (יש גם offline RL שמאמן על דאטה סטטי)."
$d.Paragraphs.Item(7).Range.Delete()

# 6) Replace the second long Hebrew body paragraph (now paragraph 8,
#    formerly paragraph 9) the same way.
$d.Paragraphs.Item(7).Range.InsertParagraphAfter()
$d.Paragraphs.Item(8).Range.Text = "באמצעות מודלי RL ניתן לאמן בין השאר רובוטים, רכבים אוטונומיים, מודלים להתמודדות עם איומי סייבר.  לאחרונה יצאו כמה שיטות אימון מודלי שפה באמצעות טכניקה שנלקחה מעולם ה-RL הנקראת RLHF. ה-ChatGPT המפורסם אומן תוך שימוש בטכניקה זו. המאמר המסוקר נשאלת השאלה האם ניתן לאמן רובוט לבצע פעולות מורכבות באמצעות מודלי שפה?  מתברר שהתשובה לשאלה הזו היא כן."
$d.Paragraphs.Item(9).Range.Delete()

# 7) Replace the third long Hebrew body paragraph (now paragraph 10,
#    formerly paragraph 11). Its new text legitimately ends with a
#    trailing space, so keeping the existing xml:space="preserve" is
#    correct - a plain text assignment is fine here.
$d.Paragraphs.Item(10).Range.Text = "המאמר לוקח מודל שפה מאומן (עם משקלים מוקפאים) ובנוסף מודל ויזואלי (מוקפא גם כן) ורותם אותם למשימת אימון זו. למשל ניתן לאמן רובוט לבצע פקודה הבאה: ״קח תפוח, בננה ולמון ותשים אותם יחד למקרר״. הגישה המוצעת היא די פשוטה. קודם כל לוקחים פקודה בשפה טבעית ובונים את השיכון (embedding) שלה באמצעות llm. בנוסף בכל שלב (נגיד אחרי כל תמונה של רובוט) מצלמים את הסביבה ומעבירים את התמונה דרך מודל ויזואלי כדי לקבל שיכון של התמונה. את ייצוג התמונה מעבירים דרך MLP מאומן(fully connected). "

# 8) Replace the fourth (last) Hebrew body paragraph (paragraph 12,
#    formerly paragraph 13). That run never had xml:space="preserve"
#    and the new text doesn't need it either, so a plain text
#    assignment is fine.
$d.Paragraphs.Item(12).Range.Text = "לאחר מכן לוקחים את ייצוג הפקודה וייצוג של כל התמונות שנבנו (אחרי ה-MLP) ומכניסים את הוקטורים האלו לאותו מודל שפה(כאילו שהם טוקנים)."

# 9) Paragraph 13 (the trailing blank "Normal" paragraph, formerly
#    paragraph 14) stays blank and unchanged. Append five new
#    paragraphs after it: blank, text, blank, text, blank-Heading2.
$p13 = $d.Paragraphs.Item($d.Paragraphs.Count)

$p13.Range.InsertParagraphAfter()
$p14 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p14.Range.Text = "ביציאה ממודל השפה מקבלים את הייצוגים ההקשריים של הטוקנים הויזואליים (תמונות). לכל טוקן ויזואלי כזה מוסיפים עוד MLP מאומן בעל שני ראשים: אחד לחישוב הפעולה הבאה והשני לחישוב פונקציית ה-value (המשערת עד כמה המצב שהרובוט נמצא בו הוא מוצלח ביחס למשימה שהוא צריך לבצע)."

$p14.Range.InsertParagraphAfter()
$p15 = $d.Paragraphs.Item($d.Paragraphs.Count)

$p15.Range.InsertParagraphAfter()
$p16 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p16.Range.Text = "בשלב האחרון מאמנים סוכן (רובוט) לבצע את הפעולות האופטימליות בהתבסס על ייצוג הפקודה ועל ייצוגי התמונות של המצבים הקודמים תוך שימוש באיזה מודיפיקציה של (PPO (proximal policy optimization הנקרא DD-PPO. פונקציית תגמול כמובן קשורה להצלחה בביצוע המשימה. כאמור מאמין שני ה-MLPs שדיברנו עליהם קודם.  נציין שבעיית RL זו היא לא פשוטה בכלל עקב מורכבות המשימה והספרסיות של התגמול (מקבלים אותו רק בסוף אחרי הרבה שלבים). למרות זאת יש תוצאות יפות."

$p16.Range.InsertParagraphAfter()
$p17 = $d.Paragraphs.Item($d.Paragraphs.Count)

$p17.Range.InsertParagraphAfter()
$p18 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p18.Style = "Heading 2"

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
